$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.09347130405523991
$ws.Range("C2").Value = 1.163400873526787
$ws.Range("D2").Value = 8.140377539043413
$ws.Range("E2").Value = 2.853134686453378
$ws.Range("F2").Value = 2.883113556558969
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.02024706560588739
$ws.Range("C3").Value = 1.092969939287157
$ws.Range("D3").Value = 5.366131528809891
$ws.Range("E3").Value = 2.31649121060493
$ws.Range("F3").Value = 2.342577598071925
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = -0.03266126848104425
$ws.Range("C4").Value = 0.9435118373712057
$ws.Range("D4").Value = 4.478804720633087
$ws.Range("E4").Value = 2.116318671805616
$ws.Range("F4").Value = 2.14053063588333
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.07925189262698597
$ws.Range("C5").Value = 1.009367654643358
$ws.Range("D5").Value = 4.81238600652716
$ws.Range("E5").Value = 2.193715115170418
$ws.Range("F5").Value = 2.218228171063641
$ws.Range("G5").Value = 43

$ws.Range("B6").Value = 0.02771576539098436
$ws.Range("C6").Value = 0.9479935622673827
$ws.Range("D6").Value = 4.345380084005096
$ws.Range("E6").Value = 2.084557527151768
$ws.Range("F6").Value = 2.109639321010291
$ws.Range("G6").Value = 42

$ws.Range("B7").Value = 0.1037329126159212
$ws.Range("C7").Value = 0.9868912652243453
$ws.Range("D7").Value = 4.448777381052547
$ws.Range("E7").Value = 2.109212502583025
$ws.Range("F7").Value = 2.132830815017939
$ws.Range("G7").Value = 41

$ws.Range("B8").Value = 0.06410150100714884
$ws.Range("C8").Value = 0.9537641320944577
$ws.Range("D8").Value = 4.474893186289131
$ws.Range("E8").Value = 2.115394333520143
$ws.Range("F8").Value = 2.141359305616774
$ws.Range("G8").Value = 40

$ws.Range("B9").Value = 0.1191595648948872
$ws.Range("C9").Value = 1.014371866898358
$ws.Range("D9").Value = 4.568013714778091
$ws.Range("E9").Value = 2.137291209633842
$ws.Range("F9").Value = 2.16186307201873
$ws.Range("G9").Value = 39

$ws.Range("B10").Value = 0.08350116669570022
$ws.Range("C10").Value = 0.9748485792014643
$ws.Range("D10").Value = 4.600399335576155
$ws.Range("E10").Value = 2.14485415251857
$ws.Range("F10").Value = 2.171997597480076
$ws.Range("G10").Value = 38

$ws.Range("B11").Value = 0.116349882878518
$ws.Range("C11").Value = 1.0065703756427
$ws.Range("D11").Value = 4.725732831804032
$ws.Range("E11").Value = 2.173875072722449
$ws.Range("F11").Value = 2.200702128189578
$ws.Range("G11").Value = 37

